$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New quest row (row 10): "the find kid quest" ---
# Write the new shared strings in the same order the target workbook has them
# appended in xl/sharedStrings.xml (H10, J10, B10, C10, then the two amended
# H6/H7 texts), so the resulting shared-string table lines up with the diff.

$ws.Range("H10").Value = "G|玛莎||告诉你，他的孩子，前几天走失了。如果你碰到了这个孩子，一定要记得把他带回来"

$ws.Range("J10").Value = "对话"
# Give J10 the same "highlight" look as the other Type cells (J6/J7 use a red
# highlight fill) but with the accent3 theme green used for this new quest.
$ws.Range("J10").Interior.Color = 255
$ws.Range("J10").Interior.ThemeColor = 7

$ws.Range("B10").Value = "走失的孩子"
$ws.Range("C10").Value = "lossboy"

# Update the two existing quest descriptions (wolf den / qiongqi) to add the
# missing trailing "|" delimiter.
$ws.Range("H6").Value = "G|塞尼斯||希望你帮助他，进入附近的|O|狼穴||并消灭狼群，如果失败了，你可以尝试反复进出本地图重试。"
$ws.Range("H7").Value = "G|塞尼斯||告诉你，再附近的森林深处，有一只神兽|O|穷奇||，如果你可以找到并击败他，会得到丰厚的回报。"

# Remaining row 10 field edits.
$ws.Range("F10").Value = 5
$ws.Range("G10").Value = 2
$ws.Range("I10").Value = 42120014
$ws.Range("Q10").Value = 100

# Reflect the selection left on the sheet after the edit.
$ws.Range("H10").Select()
